$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (column F) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 520
$wsExhibition.Range("F4").Value = 170
$wsExhibition.Range("F5").Value = 215
$wsExhibition.Range("F6").Value = 366
$wsExhibition.Range("F8").Value = 2227
$wsExhibition.Range("F9").Value = 378
$wsExhibition.Range("F10").Value = 5480

# Sheet "全部类型" (All Types) - update "想去人数" (column F) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 520
$wsAll.Range("F5").Value = 170
$wsAll.Range("F6").Value = 215
$wsAll.Range("F7").Value = 366
$wsAll.Range("F11").Value = 2227
$wsAll.Range("F12").Value = 378
$wsAll.Range("F13").Value = 5480
